$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (currently sitting in
#    the trailing empty paragraph at the end of the document).
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------
# 2. Append the new paragraphs at the end of the document:
#      <empty>
#      <empty>
#      "Intro web framework"
#      <hyperlink to MDN web frameworks article>
#      <empty>
# ---------------------------------------------------------------
$r = $d.Range($d.Content.End, $d.Content.End)
$r.Text = "`r"

$r = $d.Range($d.Content.End, $d.Content.End)
$r.Text = "`r"

$r = $d.Range($d.Content.End, $d.Content.End)
$r.Text = "Intro web framework`r"

$r = $d.Range($d.Content.End, $d.Content.End)
$r.Text = "`r"

$lastParaCount = $d.Paragraphs.Count
$hyperlinkPara = $d.Paragraphs.Item($lastParaCount - 1)
$linkRange = $d.Range($hyperlinkPara.Range.Start, $hyperlinkPara.Range.Start)
$d.Hyperlinks.Add($linkRange, "https://developer.mozilla.org/en-US/docs/Learn/Server-side/First_steps/Web_frameworks") | Out-Null

# ---------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark at the end of the
#    "Quora Lavavel vs Express" paragraph (a collapsed bookmark
#    right before the paragraph mark).
# ---------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "vs Express") {
        $target = $p
        break
    }
}

$endPos = $target.Range.End - 1
$insPos = $d.Range($endPos, $endPos)
$insPos.InsertAfter("X")
$tempRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $tempRange)
$delRange = $d.Range($endPos, $endPos + 1)
$delRange.Text = ""
